# Modified the test methods in send / request features
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename a couple of parameter tokens inside the keyword/args cell (H2)
$keywords = $ws.Range("H2").Value2
$keywords = $keywords -replace "-pphoneDescription,", "-pphoneVerificationDescription,"
$keywords = $keywords -replace "-pemailDesc,", "-pemailVerificationDescription,"
$ws.Range("H2").Value2 = $keywords

# Widen column H so the longer text still reads well
$ws.Columns.Item(8).ColumnWidth = 35

# Row 2 shrinks slightly now that the wrapped text reflows
$ws.Rows.Item(2).RowHeight = 229.5

# Active cell/selection moves from F2 to C2
$ws.Range("C2").Select()
